$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text so numeric-looking strings like "314.27"
# are stored verbatim as text (matching the source inline-string cells) instead
# of being auto-coerced into numbers by the Value setter.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.295.59'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '1.832.44'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('D5').Value = '314.27'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('E6').Value = '  +0.66%  '
$ws.Range('D7').Value = '0.4738'
$ws.Range('E7').Value = '  +1.87%  '
$ws.Range('D8').Value = '0.3683'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').Value = '0.07440'
$ws.Range('E9').Value = '  +1.15%  '
$ws.Range('D10').Value = '0.8854'
$ws.Range('E10').Value = '  +1.57%  '
$ws.Range('D11').Value = '20.48'
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').Value = '1.884.45'
$ws.Range('E12').Value = '  +5.31%  '
$ws.Range('D13').Value = '0.07322'
$ws.Range('E13').Value = '  +3.04%  '
$ws.Range('D14').Value = '5.425'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = '93.79'
$ws.Range('E15').Value = '  +2.52%  '
$ws.Range('D16').Value = '6.556'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('D18').Value = '0.000008797'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').Value = '27.601.41'
$ws.Range('E20').Value = '  +2.36%  '
$ws.Range('D21').Value = '14.77'
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('D22').Value = '5.281'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = '10.66'
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('D24').Value = '2.105.85'
$ws.Range('E24').Value = '  +3.98%  '
$ws.Range('D25').Value = '1.905'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('D26').Value = '151.83'
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('D27').Value = '18.64'
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('D28').Value = '2.139'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').Value = '117.20'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').Value = '0.08982'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('D32').Value = '0.7497'
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').Value = '1.175'
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('D34').Value = '4.540'
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('D35').Value = '2.946'
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('D36').Value = '1.011'
$ws.Range('E36').Value = '  +0.83%  '
$ws.Range('D37').Value = '1.095'
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').Value = '0.05340'
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').Value = '0.01956'
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('D40').Value = '2.420'
$ws.Range('E40').Value = '  +3.13%  '
$ws.Range('D41').Value = '2.956'
$ws.Range('E41').Value = '  -0.60%  '
$ws.Range('D42').Value = '7.238'
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').Value = '0.5290'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '0.1658'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').Value = '8.490'
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('D46').Value = '0.4915'
$ws.Range('E46').Value = '  +0.88%  '
$ws.Range('D47').Value = '10.54'
$ws.Range('E47').Value = '  +0.69%  '
$ws.Range('D48').Value = '105.13'
$ws.Range('E48').Value = '  +1.74%  '
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('D50').Value = '1.662'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').Value = '0.06296'
$ws.Range('E51').Value = '  +0.05%  '

# Restore the default (unstyled) look for column D now that the text values
# are committed, matching the original workbook's formatting.
$ws.Range("D2:D51").Style = "Normal"
